# A new weekly price observation is inserted as row 81 ("Fruta / hortaliza,
# semanal") on the "Fruta, Vega Monumental Concepción - Mango" sheet. All
# pre-existing rows from 81 down to 163 shift down by one (to 82..164); the
# sheet's used range grows from A1:T163 to A1:T164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 81, pushing the old rows 81-163 down to 82-164.
$ws.Rows("81:81").Insert()

# Populate the newly inserted row 81 with the new observation.
$ws.Range("A81").Value = 11
$ws.Range("B81").Value = "Vega Monumental Concepción"
$ws.Range("C81").Value = "Bíobío"
$ws.Range("D81").Value = 45049
$ws.Range("E81").Value = 8
$ws.Range("F81").Value = "Fruta"
$ws.Range("G81").Value = 100108
$ws.Range("H81").Value = "Tropicales y subtropicales"
$ws.Range("I81").Value = 100108002
$ws.Range("J81").Value = "Mango"
$ws.Range("K81").Value = "Sin especificar"
$ws.Range("L81").Value = "Primera"
$ws.Range("M81").Value = 110
$ws.Range("N81").Value = 7500
$ws.Range("O81").Value = 8000
$ws.Range("P81").Value = 7773
$ws.Range("Q81").Value = "$/bandeja 4 kilos"
$ws.Range("R81").Value = "Perú"
$ws.Range("S81").Value = 1943
$ws.Range("T81").Value = 4
